$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skills")

# Insert a new row above the existing CitaviR row (row 3) to host the new
# "BioMathR" R package entry, copying formatting (style + row height) from
# the row that is being pushed down so the inserted row looks identical to
# its neighbours.
$ws.Rows.Item(3).Insert()

$srcFormat = $ws.Range("A4:D4")
$dstFormat = $ws.Range("A3:D3")
$srcFormat.Copy()
$dstFormat.PasteSpecial(-4122)
$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(4).RowHeight()

# New row: BioMathR package, under the "Open Source" category
$ws.Range("A3").Value = "Open Source"
$ws.Range("B3").Value = "Open Source"
$ws.Range("C3").Value = "R package BioMathR https://schmidtpaul.github.io/BioMathR/"
$ws.Range("D3").Value = "R Paket BioMathR https://schmidtpaul.github.io/BioMathR/"

# Existing CitaviR row (now shifted to row 4): drop the trailing period
# (German column first, then English, so the shared-string table ends up
# in the same order as the source file).
$ws.Range("D4").Value = "R Paket CitaviR schmidtpaul.github.io/CitaviR/"
$ws.Range("C4").Value = "R package CitaviR https://schmidtpaul.github.io/CitaviR/"

# The author's last action was switching focus to the "Skills" tab, so make
# it the active sheet (Excel keeps the previous selection of the "Job"
# sheet at H9, but drops its "active tab" flag).
$job = $wb.Worksheets.Item("Job")
[void]$job.Range("H9").Select()

[void]$ws.Activate()
[void]$ws.Range("A1").Select()
